$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("168:168").Delete()
